$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.373.64'
$ws.Range("E2").Value = '  -0.04%  '
$ws.Range("D3").Value = '1.882.28'
$ws.Range("E3").Value = '  +0.30%  '
$ws.Range("E4").Value = '  +0.21%  '
$ws.Range("D5").Value = '0.7117'
$ws.Range("E5").Value = '  -0.06%  '
$ws.Range("D6").Value = '243.01'
$ws.Range("E6").Value = '  +0.41%  '
$ws.Range("E7").Value = '  +0.18%  '
$ws.Range("D8").Value = '0.08019'
$ws.Range("E8").Value = '  +2.91%  '
$ws.Range("D9").Value = '0.3161'
$ws.Range("E9").Value = '  +1.18%  '
$ws.Range("D10").Value = '25.09'
$ws.Range("E10").Value = '  -0.40%  '
$ws.Range("D11").Value = '0.08336'
$ws.Range("E11").Value = '  -1.40%  '
$ws.Range("D12").Value = '1.896.09'
$ws.Range("E12").Value = '  +1.28%  '
$ws.Range("D13").Value = '5.272'
$ws.Range("E13").Value = '  +0.70%  '
$ws.Range("D14").Value = '94.86'
$ws.Range("E14").Value = '  +3.94%  '
$ws.Range("E15").Value = '  +0.76%  '
$ws.Range("D16").Value = '6.365'
$ws.Range("E16").Value = '  +5.13%  '
$ws.Range("D17").Value = '0.000008668'
$ws.Range("E17").Value = '  +5.11%  '
$ws.Range("D18").Value = '29.382.00'
$ws.Range("E18").Value = '  -0.01%  '
$ws.Range("D19").Value = '243.09'
$ws.Range("E19").Value = '  +0.78%  '
$ws.Range("D20").Value = '2.151.00'
$ws.Range("E20").Value = '  +1.72%  '
$ws.Range("D21").Value = '13.34'
$ws.Range("E21").Value = '  +0.64%  '
$ws.Range("E22").Value = '  +0.24%  '
$ws.Range("D23").Value = '7.843'
$ws.Range("E23").Value = '  +0.63%  '
$ws.Range("D24").Value = '1.002'
$ws.Range("E24").Value = '  +0.16%  '
$ws.Range("D25").Value = '0.1573'
$ws.Range("E25").Value = '  -1.53%  '
$ws.Range("D26").Value = '9.101'
$ws.Range("D27").Value = '163.28'
$ws.Range("E27").Value = '  -0.04%  '
$ws.Range("E28").Value = '  +0.55%  '
$ws.Range("E29").Value = '  -0.11%  '
$ws.Range("D30").Value = '4.443'
$ws.Range("E30").Value = '  +0.29%  '
$ws.Range("D31").Value = '4.352'
$ws.Range("E31").Value = '  +0.41%  '
$ws.Range("E32").Value = '  -6.47%  '
$ws.Range("D33").Value = '0.05405'
$ws.Range("E33").Value = '  +2.14%  '
$ws.Range("D34").Value = '1.949'
$ws.Range("E34").Value = '  +0.37%  '
$ws.Range("D35").Value = '0.7752'
$ws.Range("E35").Value = '  +4.11%  '
$ws.Range("E36").Value = '  +0.60%  '
$ws.Range("E37").Value = '  -0.45%  '
$ws.Range("D38").Value = '0.01892'
$ws.Range("E38").Value = '  +1.11%  '
$ws.Range("D39").Value = '1.274.62'
$ws.Range("E39").Value = '  +4.00%  '
$ws.Range("E40").Value = '  +0.90%  '
$ws.Range("D41").Value = '6.526'
$ws.Range("E41").Value = '  +0.65%  '
$ws.Range("D42").Value = '0.9186'
$ws.Range("E42").Value = '  +2.72%  '
$ws.Range("D43").Value = '113.26'
$ws.Range("E43").Value = '  +2.64%  '
$ws.Range("D44").Value = '74.59'
$ws.Range("E44").Value = '  +2.39%  '
$ws.Range("D45").Value = '1.002'
$ws.Range("E45").Value = '  +0.18%  '
$ws.Range("D46").Value = '0.00000000129'
$ws.Range("E46").Value = '  +5.06%  '
$ws.Range("D47").Value = '2.044.19'
$ws.Range("E47").Value = '  +1.46%  '
$ws.Range("D48").Value = '1.816'
$ws.Range("E48").Value = '  -0.08%  '
$ws.Range("D49").Value = '0.5226'
$ws.Range("E49").Value = '  +0.26%  '
$ws.Range("D50").Value = '9.572'
$ws.Range("E50").Value = '  +1.81%  '
$ws.Range("D51").Value = '0.4384'
$ws.Range("E51").Value = '  +1.22%  '
